# Aggiornamento dati al 23 agosto 2021
# Append new daily rows (344-357) to the sheet, continuing the existing
# series of date / nuovi pos. / somma mobile 7gg. / somma mobile 7gg. per 100mila abitanti

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 343
$firstNewRow = $lastRow + 1

$data = @(
    @(44418, 2, 13, 131.8191036300953),
    @(44419, 1, 12, 121.6791725816264),
    @(44420, 2, 13, 131.8191036300953),
    @(44421, 0, 13, 131.8191036300953),
    @(44422, 1, 13, 131.8191036300953),
    @(44423, 0, 11, 111.5392415331576),
    @(44424, 3, 9, 91.25937943621983),
    @(44425, 0, 7, 70.9795173392821),
    @(44426, 0, 6, 60.83958629081322),
    @(44427, 1, 5, 50.69965524234435),
    @(44428, 2, 7, 70.9795173392821),
    @(44429, 1, 7, 70.9795173392821),
    @(44430, 2, 9, 91.25937943621983),
    @(44431, 4, 10, 101.3993104846887)
)

$lastNewRow = $firstNewRow + $data.Length - 1

# Copy the formatting (date style/number format/border/alignment) of the
# last existing date cell in column A down onto the new date cells.
$ws.Range("A" + $lastRow).Copy()
$ws.Range("A" + $firstNewRow + ":A" + $lastNewRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $firstNewRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
}
